$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2023" column (Q), mirroring column P's styling.
# Copy the whole P2:P6 block's formatting into Q2:Q6 first (keeps borders/
# number formats/fonts consistent with the rest of the row), then overwrite
# the values that differ for the new year.
$ws.Range("P2:P6").Copy()
$ws.Range("Q2:Q6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 3 - year header
$ws.Range("Q3").Value = 2023

# Row 4 - per-person figure (literal value, not a formula, per source data)
$ws.Range("Q4").Value = 279.01945525291825

# Row 5 - thousand tons removed
$ws.Range("Q5").Value = 1792.7

# Row 6 - average annual population
$ws.Range("Q6").Value = 6425

# Row heights adjusted per the new layout
$ws.Rows.Item(4).RowHeight = 27
$ws.Rows.Item(5).RowHeight = 27.75
